# Planilla de Métricas V2.1 Vagones de Tren
# - Corrige los datos del incremento "Tren: separarAnimalesEnVagones" (fila 22)
# - Agrega una nueva fila (23) para el incremento "Main"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Métricas")

# --- Fila 22: Tren: separarAnimalesEnVagones -> correcciones de errores ---
$ws.Range("K22").Value = 5
$ws.Range("L22").Value = 15 / (24 * 60)    # 00:15:00

# --- Fila 23: nuevo incremento "Main" ---
$ws.Range("C23").Value = "Main"
$ws.Range("F23").Value = 10
$ws.Range("G23").Value = 2 / (24 * 60)                    # 00:02:00 (Tiempo estimado)
$ws.Range("H23").Value = (17 + 25 / 60) / 24               # 17:25 Hora Inicio
$ws.Range("I23").Value = (17 + 26 / 60) / 24                # 17:26 Hora Fin
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("M23").Value = 7

# --- Selección activa tal cual queda en el archivo final ---
$ws.Range("M23").Select()
